$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The template row (row 4) already carries the correct cell styles for
# columns A:K used by the data rows; copy its formatting down into the
# two new/blank rows (5 and 6) before filling in values so the new rows
# pick up the same styles as the rest of the table. Row 5 only needs
# columns A:J (no trailing K cell), row 6 needs the full A:K range.
$ws.Range("A4:J4").Copy()
$ws.Range("A5:J5").PasteSpecial(-4122)
$ws.Range("A4:K4").Copy()
$ws.Range("A6:K6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 5: Crumpet GEF
$ws.Range("A5").Value = "Crumpet GEF"
$ws.Range("B5").Value = 20001371
$ws.Range("C5").Value = "Crumpet exporter"
$ws.Range("D5").Value = "GBP"
$ws.Range("E5").Value = 7000000
$ws.Range("F5").Value = 3938753.8
$ws.Range("G5").Value = 777
$ws.Range("H5").Value = 456
$ws.Range("I5").Value = "GBP"
$ws.Range("J5").Value = "GBP"

# Row 6: Scone GEF
$ws.Range("A6").Value = "Scone GEF"
$ws.Range("B6").Value = 20001371
$ws.Range("C6").Value = "Scone exporter"
$ws.Range("D6").Value = "GBP"
$ws.Range("E6").Value = 770000
$ws.Range("F6").Value = 761579.37
$ws.Range("G6").Value = 777
$ws.Range("H6").Value = 456.77
$ws.Range("I6").Value = "GBP"
$ws.Range("J6").Value = "GBP"

# Update selection to match the target state
$ws.Range("A5:J6").Select()
